# "Requests should be completed" - fill in the outstanding request-form
# fields (request number, date, incoming numbers, test duration) and tidy
# up the now-redundant placeholder rows in both test tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header: request number + corrected date
# (the merged A6:F6 / A15:F15 ranges carry the text on every cell in the
# underlying file, not just the merge anchor, so mirror that here too)
# ---------------------------------------------------------------------
$ws.Range("A6:F6").Value = "ЗАЯВКА № A/B50278 / Дата 07.05.2016"
$ws.Range("A15:F15").Value = "ЗАЯВКА № A/B50278 / Дата 07.05.2016"

# Fix the double space in the "Вид на пробата" column header (both tables)
$ws.Range("B9").Value = "Вид на пробата"
$ws.Range("B18").Value = "Вид на пробата"

# ---------------------------------------------------------------------
# 2. Table 1 (rows 9-13): fill in the incoming number and collapse the
#    now-unused placeholder / duplicate rows.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "A/B50278"
$ws.Range("B10").Value = "1. p1"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "234"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3234"
$ws.Range("E10").ClearContents()
$ws.Range("F10").ClearContents()

$ws.Range("A11").ClearContents()
$ws.Range("B11").Value = "2. p2"
$ws.Range("C11").Value = "E. Coli"
$ws.Range("D11").Value = "ISO 16649-2"
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()

# Row 12 duplicated row 11's data and is no longer needed.
$ws.Range("A12:F12").ClearContents()

$ws.Range("A13:C13").Value = "Срок за изпитване: 11 дни"

# ---------------------------------------------------------------------
# 3. Table 2 (rows 17-21): same treatment - fill incoming number, move
#    the sample row up into the previously-empty row 19, drop old row 20.
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "A/B50278"
$ws.Range("B19").Value = "1. p1"
$ws.Range("C19").Value = "Бацилус"
$ws.Range("D19").Value = " БДС ..."
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()

$ws.Range("A20:F20").ClearContents()

$ws.Range("A21:C21").Value = "Срок за изпитване: 11 дни"

# ---------------------------------------------------------------------
# 4. Formatting to match the "completed" look of the form.
# ---------------------------------------------------------------------

# A6 / A15 - bold 14pt Times New Roman, centered (request line)
foreach ($addr in @("A6", "A15")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Times New Roman"
    $c.Font.Size = 14
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
}

# A8 / A17 - bold 14pt Times New Roman, indented section titles
foreach ($addr in @("A8", "A17")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Times New Roman"
    $c.Font.Size = 14
    $c.Font.Bold = $true
    $c.IndentLevel = 2
}

# Row 9 / Row 18 header rows - 12pt Times New Roman, boxed, centered
foreach ($rowAddr in @("A9:F9", "A18:F18")) {
    $c = $ws.Range($rowAddr)
    $c.Font.Name = "Times New Roman"
    $c.Font.Size = 12
    $c.Font.Bold = $false
    $c.Borders.LineStyle = 1
    $c.Borders.Weight = 2
    $c.HorizontalAlignment = -4108
}

# Data rows (10, 11, 19): A & F columns get a thin box + wrap + vcenter,
# B:E columns get boxed, centered, wrapped 12pt Times New Roman text.
foreach ($rowNum in 10, 11, 19) {
    $edge = $ws.Range("A" + $rowNum + ":F" + $rowNum)
    $edge.Borders.LineStyle = 1
    $edge.Borders.Weight = 2
    $edge.WrapText = $true
    $edge.VerticalAlignment = -4108

    $inner = $ws.Range("B" + $rowNum + ":E" + $rowNum)
    $inner.Font.Name = "Times New Roman"
    $inner.Font.Size = 12
    $inner.HorizontalAlignment = -4108
}

# A13 / D13 / A21 / D21 - indented plain text
foreach ($addr in @("A13", "D13", "A21", "D21")) {
    $c = $ws.Range($addr)
    $c.IndentLevel = 3
}
